# Update the "Handläggningsdatum" (column C) date value for rows 2-23
# from 2023-09-19 (serial 45188) to 2023-09-20 (serial 45189).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
